# Apply the "new lm" metrics update to the results sheet.
# - Column A (model ids) gets re-ordered: row 26 (model_8_5_24) keeps its
#   place, every other row is remapped to a different model id.
# - Columns B..Q (the metric values) become identical across every data
#   row (2..26), taking on the values that used to belong to model_8_5_4
#   (originally row 2) combined with the new model_8_5_24-style tail
#   metrics for row 26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: reordered model identifiers ---------------------------------
$ws.Range("A2").Value  = "model_8_5_0"
$ws.Range("A3").Value  = "model_8_5_22"
$ws.Range("A4").Value  = "model_8_5_21"
$ws.Range("A5").Value  = "model_8_5_20"
$ws.Range("A6").Value  = "model_8_5_19"
$ws.Range("A7").Value  = "model_8_5_18"
$ws.Range("A8").Value  = "model_8_5_17"
$ws.Range("A9").Value  = "model_8_5_16"
$ws.Range("A10").Value = "model_8_5_15"
$ws.Range("A11").Value = "model_8_5_14"
$ws.Range("A12").Value = "model_8_5_13"
$ws.Range("A13").Value = "model_8_5_23"
$ws.Range("A14").Value = "model_8_5_12"
$ws.Range("A15").Value = "model_8_5_10"
$ws.Range("A16").Value = "model_8_5_9"
$ws.Range("A17").Value = "model_8_5_8"
$ws.Range("A18").Value = "model_8_5_7"
$ws.Range("A19").Value = "model_8_5_6"
$ws.Range("A20").Value = "model_8_5_5"
$ws.Range("A21").Value = "model_8_5_4"
$ws.Range("A22").Value = "model_8_5_3"
$ws.Range("A23").Value = "model_8_5_2"
$ws.Range("A24").Value = "model_8_5_1"
$ws.Range("A25").Value = "model_8_5_11"
$ws.Range("A26").Value = "model_8_5_24"

# --- Columns B..Q: identical metric values for every data row --------------
# (NOTE: the PowerShell parser here doesn't accept bare scientific-notation
#  literals like `9.6e-07`, so those are built via [double]"...")
$ws.Range("B2:B26").Value  = 0.9999989698596451
$ws.Range("C2:C26").Value  = 0.9990244665669379
$ws.Range("D2:D26").Value  = 0.9999662338913162
$ws.Range("E2:E26").Value  = 0.9999996448331083
$ws.Range("F2:F26").Value  = 0.9999968605020622
$ws.Range("G2:G26").Value  = [double]"9.615908033740045e-07"
$ws.Range("H2:H26").Value  = 0.0009106176387947784
$ws.Range("I2:I26").Value  = [double]"4.194864220611975e-06"
$ws.Range("J2:J26").Value  = [double]"3.3416215023955e-07"
$ws.Range("K2:K26").Value  = [double]"2.264513185425763e-06"
$ws.Range("L2:L26").Value  = [double]"5.130350677945762e-05"
$ws.Range("M2:M26").Value  = 0.0009806073645318011
$ws.Range("N2:N26").Value  = 1.000024723368518
$ws.Range("O2:O26").Value  = 0.001022353850177004
$ws.Range("P2:P26").Value  = 77.70935367418031
$ws.Range("Q2:Q26").Value  = 108.1812492958853
